# Third commit: after the "changing the file now 1.1" paragraph, add a
# blank paragraph and a new paragraph with "new changes of third commit".
# The paragraph mark of the edited / inserted paragraphs picks up the
# Calibri / en-language run formatting already used in the document, as
# Word does when you type at the end of a formatted run and press Enter.

$d = $word.ActiveDocument

# Locate the paragraph that holds "changing the file now 1.1" robustly
# (rather than hard-coding an index).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*changing the file now 1.1*") {
        $target = $p
    }
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:lang w:val="en"/></w:rPr>'

$newXml =
    '<w:p ' + $ns + '><w:pPr>' + $rPr + '</w:pPr>' +
        '<w:r>' + $rPr + '<w:t>changing the file now 1.1</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $ns + '><w:pPr>' + $rPr + '</w:pPr></w:p>' +
    '<w:p ' + $ns + '><w:r>' + $rPr + '<w:t>new changes of third commit</w:t></w:r></w:p>'

$target.Range.InsertXML($newXml)
